# step_2 updated and reran
# The DiSCoVER "top drugs (cerebellar stem cell control)" table slide
# (slide 2) is re-run/duplicated and appended as a new slide at the end
# of the deck.

$p = $ppt.ActivePresentation

# Slide 2 is "DiSCoVER: top drugs (cerebellar stem cell control)" -
# duplicate it and move the duplicate to the end of the slide list.
$src = $p.Slides.Item(2)
$dup = $src.Duplicate()
$newSlide = $dup.Item(1)
$newSlide.MoveTo($p.Slides.Count)
